$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 108-114 in place, and append new rows 115-118
# with the final post-edit values (data table row realignment).

# Row 108
$ws.Range("A108").Value = 1
$ws.Range("B108").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C108").Value = "Arica y Parinacota"
$ws.Range("D108").Value = 44918
$ws.Range("E108").Value = 15
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100109
$ws.Range("H108").Value = "Uva"
$ws.Range("I108").Value = 100109001
$ws.Range("J108").Value = "Uva"
$ws.Range("K108").Value = "Superior Seedless"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 200
$ws.Range("N108").Value = 20000
$ws.Range("O108").Value = 20000
$ws.Range("P108").Value = 20000
$ws.Range("Q108").Value = "$/caja 12 kilos"
$ws.Range("R108").Value = "Región de Coquimbo"
$ws.Range("S108").Value = 1667
$ws.Range("T108").Value = 12

# Row 109
$ws.Range("A109").Value = 1
$ws.Range("B109").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C109").Value = "Arica y Parinacota"
$ws.Range("D109").Value = 44918
$ws.Range("E109").Value = 15
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100109
$ws.Range("H109").Value = "Uva"
$ws.Range("I109").Value = 100109001
$ws.Range("J109").Value = "Uva"
$ws.Range("K109").Value = "Superior Seedless"
$ws.Range("L109").Value = "Segunda"
$ws.Range("M109").Value = 350
$ws.Range("N109").Value = 22000
$ws.Range("O109").Value = 22000
$ws.Range("P109").Value = 22000
$ws.Range("Q109").Value = "$/caja 12 kilos"
$ws.Range("R109").Value = "Región de Coquimbo"
$ws.Range("S109").Value = 1833
$ws.Range("T109").Value = 12

# Row 110
$ws.Range("A110").Value = 1
$ws.Range("B110").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C110").Value = "Arica y Parinacota"
$ws.Range("D110").Value = 44307
$ws.Range("E110").Value = 15
$ws.Range("F110").Value = "Fruta"
$ws.Range("G110").Value = 100109
$ws.Range("H110").Value = "Uva"
$ws.Range("I110").Value = 100109001
$ws.Range("J110").Value = "Uva"
$ws.Range("K110").Value = "Rosada pastilla"
$ws.Range("L110").Value = "Segunda"
$ws.Range("M110").Value = 260
$ws.Range("N110").Value = 16000
$ws.Range("O110").Value = 17000
$ws.Range("P110").Value = 16500
$ws.Range("Q110").Value = "$/bandeja 18 kilos"
$ws.Range("R110").Value = "Región de Coquimbo"
$ws.Range("S110").Value = 917
$ws.Range("T110").Value = 18

# Row 111
$ws.Range("A111").Value = 1
$ws.Range("B111").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C111").Value = "Arica y Parinacota"
$ws.Range("D111").Value = 44665
$ws.Range("E111").Value = 15
$ws.Range("F111").Value = "Fruta"
$ws.Range("G111").Value = 100109
$ws.Range("H111").Value = "Uva"
$ws.Range("I111").Value = 100109001
$ws.Range("J111").Value = "Uva"
$ws.Range("K111").Value = "Rosada pastilla"
$ws.Range("L111").Value = "Primera"
$ws.Range("M111").Value = 200
$ws.Range("N111").Value = 17000
$ws.Range("O111").Value = 18000
$ws.Range("P111").Value = 17500
$ws.Range("Q111").Value = "$/bandeja 18 kilos"
$ws.Range("R111").Value = "Región de Coquimbo"
$ws.Range("S111").Value = 972
$ws.Range("T111").Value = 18

# Row 112
$ws.Range("A112").Value = 1
$ws.Range("B112").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C112").Value = "Arica y Parinacota"
$ws.Range("D112").Value = 44665
$ws.Range("E112").Value = 15
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100109
$ws.Range("H112").Value = "Uva"
$ws.Range("I112").Value = 100109001
$ws.Range("J112").Value = "Uva"
$ws.Range("K112").Value = "Thompson seedless"
$ws.Range("L112").Value = "Primera"
$ws.Range("M112").Value = 300
$ws.Range("N112").Value = 17000
$ws.Range("O112").Value = 18000
$ws.Range("P112").Value = 17500
$ws.Range("Q112").Value = "$/bandeja 18 kilos"
$ws.Range("R112").Value = "Región de Coquimbo"
$ws.Range("S112").Value = 972
$ws.Range("T112").Value = 18

# Row 113
$ws.Range("A113").Value = 1
$ws.Range("B113").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C113").Value = "Arica y Parinacota"
$ws.Range("D113").Value = 44545
$ws.Range("E113").Value = 15
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100109
$ws.Range("H113").Value = "Uva"
$ws.Range("I113").Value = 100109001
$ws.Range("J113").Value = "Uva"
$ws.Range("K113").Value = "Superior Seedless"
$ws.Range("L113").Value = "Primera"
$ws.Range("M113").Value = 200
$ws.Range("N113").Value = 17000
$ws.Range("O113").Value = 18000
$ws.Range("P113").Value = 17500
$ws.Range("Q113").Value = "$/bandeja 10 kilos"
$ws.Range("R113").Value = "Región de O'Higgins"
$ws.Range("S113").Value = 1750
$ws.Range("T113").Value = 10

# Row 114
$ws.Range("A114").Value = 1
$ws.Range("B114").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C114").Value = "Arica y Parinacota"
$ws.Range("D114").Value = 44679
$ws.Range("E114").Value = 15
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100109
$ws.Range("H114").Value = "Uva"
$ws.Range("I114").Value = 100109001
$ws.Range("J114").Value = "Uva"
$ws.Range("K114").Value = "Rosada pastilla"
$ws.Range("L114").Value = "Segunda"
$ws.Range("M114").Value = 250
$ws.Range("N114").Value = 19000
$ws.Range("O114").Value = 20000
$ws.Range("P114").Value = 19500
$ws.Range("Q114").Value = "$/bandeja 18 kilos"
$ws.Range("R114").Value = "Región de Coquimbo"
$ws.Range("S114").Value = 1083
$ws.Range("T114").Value = 18

# Row 115
$ws.Range("A115").Value = 1
$ws.Range("B115").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C115").Value = "Arica y Parinacota"
$ws.Range("D115").Value = 44679
$ws.Range("E115").Value = 15
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100109
$ws.Range("H115").Value = "Uva"
$ws.Range("I115").Value = 100109001
$ws.Range("J115").Value = "Uva"
$ws.Range("K115").Value = "Thompson seedless"
$ws.Range("L115").Value = "Primera"
$ws.Range("M115").Value = 270
$ws.Range("N115").Value = 19000
$ws.Range("O115").Value = 20000
$ws.Range("P115").Value = 19500
$ws.Range("Q115").Value = "$/bandeja 18 kilos"
$ws.Range("R115").Value = "Región de Coquimbo"
$ws.Range("S115").Value = 1083
$ws.Range("T115").Value = 18

# Row 116
$ws.Range("A116").Value = 1
$ws.Range("B116").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C116").Value = "Arica y Parinacota"
$ws.Range("D116").Value = 44245
$ws.Range("E116").Value = 15
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100109
$ws.Range("H116").Value = "Uva"
$ws.Range("I116").Value = 100109001
$ws.Range("J116").Value = "Uva"
$ws.Range("K116").Value = "Rosada pastilla"
$ws.Range("L116").Value = "Segunda"
$ws.Range("M116").Value = 250
$ws.Range("N116").Value = 17000
$ws.Range("O116").Value = 18000
$ws.Range("P116").Value = 17500
$ws.Range("Q116").Value = "$/bandeja 18 kilos"
$ws.Range("R116").Value = "Región de Coquimbo"
$ws.Range("S116").Value = 972
$ws.Range("T116").Value = 18

# Row 117
$ws.Range("A117").Value = 1
$ws.Range("B117").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C117").Value = "Arica y Parinacota"
$ws.Range("D117").Value = 44650
$ws.Range("E117").Value = 15
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100109
$ws.Range("H117").Value = "Uva"
$ws.Range("I117").Value = 100109001
$ws.Range("J117").Value = "Uva"
$ws.Range("K117").Value = "Rosada pastilla"
$ws.Range("L117").Value = "Segunda"
$ws.Range("M117").Value = 270
$ws.Range("N117").Value = 19000
$ws.Range("O117").Value = 20000
$ws.Range("P117").Value = 19500
$ws.Range("Q117").Value = "$/bandeja 18 kilos"
$ws.Range("R117").Value = "Región de Coquimbo"
$ws.Range("S117").Value = 1083
$ws.Range("T117").Value = 18

# Row 118
$ws.Range("A118").Value = 1
$ws.Range("B118").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C118").Value = "Arica y Parinacota"
$ws.Range("D118").Value = 44650
$ws.Range("E118").Value = 15
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100109
$ws.Range("H118").Value = "Uva"
$ws.Range("I118").Value = 100109001
$ws.Range("J118").Value = "Uva"
$ws.Range("K118").Value = "Thompson seedless"
$ws.Range("L118").Value = "Segunda"
$ws.Range("M118").Value = 300
$ws.Range("N118").Value = 16000
$ws.Range("O118").Value = 17000
$ws.Range("P118").Value = 16500
$ws.Range("Q118").Value = "$/bandeja 18 kilos"
$ws.Range("R118").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S118").Value = 917
$ws.Range("T118").Value = 18

# Ensure date-style formatting (style index 2 / numFmt 165) carries to the newly appended rows
$ws.Range("D116:D118").NumberFormat = $ws.Range("D114").NumberFormat

